$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ G=19.750268; H=59.250804; I=0.2138888518073023; J=0.2138888518073023; K=3; L=1; M=0.07585433333333333; N=0.227563; O=0.0240359804038997; P=0.0240359804038997; Q=1.498143412294666; R=13.483290710652; S=0.005141028250652926; T=0.005141028250652925 }
    3  = @{ G=19.750268; H=59.250804; I=0.2138888518073023; J=0.2138888518073023; O=0.9352348927362568; P=0.9352348927362568; Q=58.292441995568; R=524.631977960112; S=0.2000363173774835; T=0.2000363173774835 }
    4  = @{ G=19.750268; H=59.250804; I=0.2138888518073023; J=0.2138888518073023; M=0.1285356666666667; N=0.385607; O=0.04072912685984344; P=0.04072912685984344; Q=2.538613864225333; R=22.847524778028; S=0.00871150617916587; T=0.00871150617916587 }
    5  = @{ I=0.6395228081370402; J=0.6395228081370402; K=3; L=1; M=0.07585433333333333; N=0.227563; O=0.0240359804038997; P=0.0240359804038997; Q=4.479414770461555; R=40.31473293415399; S=0.01537155768422881; T=0.0153715576842288 }
    6  = @{ I=0.6395228081370402; J=0.6395228081370402; O=0.9352348927362568; P=0.9352348927362568; S=0.5981040448704346; T=0.5981040448704346 }
    7  = @{ I=0.6395228081370402; J=0.6395228081370402; M=0.1285356666666667; N=0.385607; O=0.04072912685984344; P=0.04072912685984344; Q=7.590397786078445; R=68.31358007470601; S=0.02604720558237683; T=0.02604720558237682 }
    8  = @{ G=13.53581066666667; H=40.607432; I=0.1465883400556574; J=0.1465883400556574; K=3; L=1; M=0.07585433333333333; N=0.227563; O=0.0240359804038997; P=0.0240359804038997; Q=1.026749894246222; R=9.240749048215999; S=0.003523394469017968; T=0.003523394469017967 }
    9  = @{ G=13.53581066666667; H=40.607432; I=0.1465883400556574; J=0.1465883400556574; O=0.9352348927362568; P=0.9352348927362568; Q=39.95062032321067; R=359.555582908896; S=0.1370945304883387; T=0.1370945304883387 }
    10 = @{ G=13.53581066666667; H=40.607432; I=0.1465883400556574; J=0.1465883400556574; M=0.1285356666666667; N=0.385607; O=0.04072912685984344; P=0.04072912685984344; Q=1.739834447913778; R=15.658510031224; S=0.005970415098300741; T=0.005970415098300741 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
